$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 11.0
$ws.Range("B12").Value = "Tuesday, Jan 10"
$ws.Range("C12").Value = "1:00 PM"
$ws.Range("D12").Value = "FR9885"
$ws.Range("E12").Value = "Brussels"
$ws.Range("F12").Value = "(CRL)"
$ws.Range("G12").Value = "Ryanair "
$ws.Range("H12").Value = "B738"
$ws.Range("I12").Value = "(EI-DHN)"
$ws.Range("J12").Value = "1:15 PM"
$ws.Range("K12").Font.Size = 11
$ws.Range("L12").Value = "0 hours, 15 minutes"
$ws.Range("M12").Font.Size = 11
